$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "ee"
$ws.Range("B6").Value = 200
$ws.Range("A7").Value = "ff"
$ws.Range("B7").Value = 100

$ws.Range("B7").Select()
